# Update "想去人数" (F2/F3) values on the 展览 (Exhibition) and 全部类型 (All Types)
# sheets to reflect newly generated output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 727
    $ws.Range("F3").Value = 4093
}
